$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WMT_Extract_SA")
$ws.Range("A2").Value = "D42237"
Write-Host "done"
